$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.647.55'
$ws.Range("E2").Value = '  +3.74%  '
$ws.Range("D3").Value = '3.503.92'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.27%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.502.89'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.589'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +4.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.26%  '
$ws.Range("D13").Value = '4.111.82'
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.18%  '
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("D17").Value = '66.677.66'
$ws.Range("E17").Value = '  +3.71%  '
$ws.Range("D18").Value = '3.511.07'
$ws.Range("E18").Value = '  +2.35%  '
$ws.Range("E19").Value = '  +3.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +3.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.50%  '
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +6.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").Value = '  +2.98%  '
$ws.Range("E39").Value = '  +4.80%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.63%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0749'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.84%  '
$ws.Range("D44").Value = '2.817.16'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '353.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.13%  '
